# Updated cryptos list on Mon Dec  4 17:20:14 UTC 2023 with GitHub Actions
# Refresh per-coin Price (D) and Volume(1h) (E) columns, and the two pairs
# of rows whose ranking order swapped (FTXToken/VeChain and
# TrustWalletToken/Aave), to match the latest coinranking.com scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.686.06'
$ws.Range("E2").Value = '  +5.34%  '
$ws.Range("D3").Value = '2.225.07'
$ws.Range("E3").Value = '  +3.27%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.88'
$ws.Range("E5").Value = '  +1.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.11'
$ws.Range("E7").Value = '  -2.56%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +3.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.68'
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0887'
$ws.Range("E11").Value = '  +5.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.103'
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").Value = '2.555.47'
$ws.Range("E13").Value = '  +3.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.67'
$ws.Range("E14").Value = '  -1.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.71'
$ws.Range("E15").Value = '  +0.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.798'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.56'
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("D18").Value = '2.226.81'
$ws.Range("E18").Value = '  +3.34%  '
$ws.Range("D19").Value = '41.584.66'
$ws.Range("E19").Value = '  +5.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.88'
$ws.Range("E20").Value = '  +1.73%  '
$ws.Range("D21").Value = '0.0₃0894'
$ws.Range("E21").Value = '  +5.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.05'
$ws.Range("E22").Value = '  +0.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.71'
$ws.Range("E23").Value = '  +9.82%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").Value = '  +2.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.33'
$ws.Range("E26").Value = '  +0.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.57'
$ws.Range("E27").Value = '  +1.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.21'
$ws.Range("E28").Value = '  -1.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.140'
$ws.Range("E29").Value = '  +1.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.94'
$ws.Range("E30").Value = '  +1.81%  '
$ws.Range("E31").Value = '  -0.59%  '
$ws.Range("E32").Value = '  -2.48%  '
$ws.Range("E33").Value = '  +0.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.95'
$ws.Range("E34").Value = '  +5.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.62'
$ws.Range("E35").Value = '  +1.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0623'
$ws.Range("E36").Value = '  +1.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.57'
$ws.Range("E37").Value = '  -5.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.68'
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("E39").Value = '  -1.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.000247'
$ws.Range("E40").Value = '  +29.41%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("B42").Value = 'FTXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.88'
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0239'
$ws.Range("E43").Value = '  +5.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.60'
$ws.Range("E44").Value = '  +9.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0981'
$ws.Range("E45").Value = '  +6.83%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.22'
$ws.Range("E46").Value = '  +1.67%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '99.03'
$ws.Range("E47").Value = '  -3.39%  '
$ws.Range("D48").Value = '1.463.90'
$ws.Range("E48").Value = '  -3.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.48'
$ws.Range("E49").Value = '  -6.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.80'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("E51").Value = '  -1.01%  '
